# Weekly price update: a new record (date serial 44847, i.e. 2022-10-13) is
# added at the top of the date-ordered data table (row 9), pushing every
# existing record from row 9 downward by one row (to rows 10-76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 9. This shifts the existing rows 9-75
# (and their formatting, e.g. the date style on column D) down to rows 10-76.
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value  = 1
$ws.Cells.Item(9, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value  = 44847
$ws.Cells.Item(9, 5).Value  = 15
$ws.Cells.Item(9, 6).Value  = 100112009
$ws.Cells.Item(9, 7).Value  = "Acelga"
$ws.Cells.Item(9, 8).Value  = "Sin especificar"
$ws.Cells.Item(9, 9).Value  = "Primera"
$ws.Cells.Item(9, 10).Value = 250
$ws.Cells.Item(9, 11).Value = 1400
$ws.Cells.Item(9, 12).Value = 1500
$ws.Cells.Item(9, 13).Value = 1450
$ws.Cells.Item(9, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(9, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 16).Value = 483
$ws.Cells.Item(9, 17).Value = 3
$ws.Cells.Item(9, 18).Value = "Hortaliza"
